# RandomFloat66_QuickSortTimes.csv.xlsx — fix the sorting timing figures for
# the first two rows (5000- and 10000-element runs): Avg_Time_ms in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 0.58100923
$ws.Range("D3").Value = 1.3228258

# Recalculate the workbook so every dependent (chart caches included) is
# refreshed against the corrected figures.
$excel.CalculateFullRebuild()
